# The workbook gained one new weekly price record for Cilantro
# (Femacal de La Calera) that needs to be inserted right above the
# existing row 321, pushing every subsequent row (321-401) down by one
# (so the sheet grows from A1:R401 to A1:R402).
#
# Inserting a whole row preserves all the other, unrelated rows and
# lets Excel manage the dimension/row-shift bookkeeping automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 321:401 down to 322:402, leaving a blank row 321 behind.
$ws.Rows("321:321").Insert()

# Fill the newly freed row 321 with the new record's data.
$ws.Range("A321").Value = 3
$ws.Range("B321").Value = "Femacal de La Calera"
$ws.Range("C321").Value = "Coquimbo"
$ws.Range("D321").Value = 44782
$ws.Range("E321").Value = 5
$ws.Range("F321").Value = 100112040
$ws.Range("G321").Value = "Cilantro"
$ws.Range("H321").Value = "Sin especificar"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 210
$ws.Range("K321").Value = 4500
$ws.Range("L321").Value = 5000
$ws.Range("M321").Value = 4762
$ws.Range("N321").Value = "$/docena de atados (3 kilos)"
$ws.Range("O321").Value = "Provincia de Quillota"
$ws.Range("P321").Value = 1587
$ws.Range("Q321").Value = 3
$ws.Range("R321").Value = "Hortaliza"
